# Update the two non-gaussian ranova tables:
#  1. Header row height grows from 571 -> 637 twips (28.55pt -> 31.85pt)
#     to accommodate the re-run results (1000 iterations vs 10).
#  2. The chi-squared symbol in the header ("χ") is re-saved using the
#     mis-encoded byte sequence that appears in the regenerated output
#     ("Ï‡" = UTF-8 bytes of χ misread as Windows-1252).

$d = $word.ActiveDocument

# --- 1. Row height fix -----------------------------------------------
# Twips -> points: Word's Row.Height property is expressed in points.
$targetHeightTwips = 637
$targetHeightPts = $targetHeightTwips / 20.0

$sourceHeightPts = 571 / 20.0
for ($i = 1; $i -le $d.Tables.Count; $i++) {
    $tbl = $d.Tables.Item($i)
    $headerRow = $tbl.Rows.Item(1)
    if ([math]::Abs($headerRow.Height - $sourceHeightPts) -lt 0.01) {
        $headerRow.Height = $targetHeightPts
    }
}

# --- 2. χ -> Ï‡ text fix -----------------------------------------------
$oldChar = [string]([char]0x03C7)
$newChars = [string]([char]0x00CF) + [string]([char]0x2021)

$d.Content.Find.Execute($oldChar, $true, $false, $false, $false, $false, `
                         $true, 1, $false, $newChars, 2) | Out-Null
